$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Jengibre" table.
# It belongs chronologically before the existing row 170, so we insert a
# new row there (which pushes every following row down by one, and the
# former last row ends up duplicated correctly at the new final row).
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record's data.
$ws.Cells.Item(170, 1).Value = 10
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170, 3).Value = "La Araucanía"
$ws.Cells.Item(170, 4).Value = 45090
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 100114007
$ws.Cells.Item(170, 7).Value = "Jengibre"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 20
$ws.Cells.Item(170, 11).Value = 24000
$ws.Cells.Item(170, 12).Value = 24000
$ws.Cells.Item(170, 13).Value = 24000
$ws.Cells.Item(170, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(170, 15).Value = "Perú"
$ws.Cells.Item(170, 16).Value = 1846
$ws.Cells.Item(170, 17).Value = 13
$ws.Cells.Item(170, 18).Value = "Hortaliza"
